# Restore C10 on the "Rules" sheet of the Decision Table sample from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
